$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $ref, $val) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $val
}

Set-TextValue $ws "D2" "322.85"
Set-TextValue $ws "E2" "-2.20%"
Set-TextValue $ws "D3" "39.71"
Set-TextValue $ws "E3" "-0.84%"
Set-TextValue $ws "D4" "5.879"
Set-TextValue $ws "E4" "11.39%"
Set-TextValue $ws "D5" "0.08019"
Set-TextValue $ws "E5" "-0.94%"
Set-TextValue $ws "D6" "4.587"
Set-TextValue $ws "E6" "1.42%"
Set-TextValue $ws "D7" "8.658"
Set-TextValue $ws "E7" "0.17%"
Set-TextValue $ws "D8" "1.929"
Set-TextValue $ws "E8" "0.75%"
Set-TextValue $ws "E9" "-0.43%"
Set-TextValue $ws "D10" "0.1273"
Set-TextValue $ws "E10" "-8.17%"
Set-TextValue $ws "D11" "0.1968"
Set-TextValue $ws "E11" "0.47%"
Set-TextValue $ws "D12" "8.726"
Set-TextValue $ws "E12" "20.01%"
Set-TextValue $ws "D13" "0.09217"
Set-TextValue $ws "E13" "0.61%"
Set-TextValue $ws "D14" "0.03535"
Set-TextValue $ws "E14" "3.06%"
Set-TextValue $ws "E15" "0.05%"
Set-TextValue $ws "D16" "0.001297"
Set-TextValue $ws "E16" "-7.50%"
Set-TextValue $ws "D17" "0.006120"
Set-TextValue $ws "E17" "1.05%"
Set-TextValue $ws "D18" "3.348"
Set-TextValue $ws "E18" "-0.36%"
Set-TextValue $ws "D19" "2.943"
Set-TextValue $ws "E19" "-0.53%"
Set-TextValue $ws "E20" "1.07%"
Set-TextValue $ws "D21" "0.1418"
Set-TextValue $ws "E21" "7.95%"
Set-TextValue $ws "E22" "-6.23%"
Set-TextValue $ws "D23" "0.04413"
Set-TextValue $ws "E23" "-0.71%"
Set-TextValue $ws "D24" "0.001261"
Set-TextValue $ws "E24" "3.18%"
Set-TextValue $ws "D25" "0.004371"
Set-TextValue $ws "E25" "0.30%"
Set-TextValue $ws "D26" "0.0001141"
Set-TextValue $ws "E26" "-11.63%"
Set-TextValue $ws "D39" "0.02448"
Set-TextValue $ws "E39" "-2.91%"
Set-TextValue $ws "E40" "0.39%"
Set-TextValue $ws "D41" "0.007420"
Set-TextValue $ws "E41" "-3.97%"
Set-TextValue $ws "D42" "0.009471"
Set-TextValue $ws "E42" "6.19%"
Set-TextValue $ws "D43" "0.1406"
Set-TextValue $ws "E43" "-1.67%"
Set-TextValue $ws "D44" "0.002122"
Set-TextValue $ws "E44" "-2.31%"
Set-TextValue $ws "D45" "0.009867"
Set-TextValue $ws "E45" "9.68%"
Set-TextValue $ws "D46" "0.00006738"
Set-TextValue $ws "E46" "1.97%"
Set-TextValue $ws "E47" "0.02%"
Set-TextValue $ws "D48" "0.003001"
Set-TextValue $ws "E48" "-10.23%"
Set-TextValue $ws "D50" "0.00002100"
Set-TextValue $ws "E50" "0.02%"
Set-TextValue $ws "D51" "0.0002000"
Set-TextValue $ws "E51" "0.02%"

Write-Host "Applied cryptos price/volume updates"
